$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quiz answer/color data cells
$ws.Range("D3").Value = "right"
$ws.Range("E3").Value = "orange"

$ws.Range("D10").Value = "right"
$ws.Range("E10").Value = "green"

$ws.Range("B14").Value = "Charcoal"
$ws.Range("C14").Value = "Flood"
$ws.Range("D14").Value = "left"
$ws.Range("E14").Value = "#ff3131"

$ws.Range("D18").Value = "left"
$ws.Range("E18").Value = "pink"

$ws.Range("E19").Value = "yellow"

$ws.Range("D21").Value = "right"
$ws.Range("E21").Value = "#ff3131"

$ws.Range("D29").Value = "right"
$ws.Range("E29").Value = "green"

$ws.Range("D31").Value = "left"
$ws.Range("E31").Value = "yellow"

# Update the selected cell / scroll position as reflected in the saved view
$ws.Range("H10").Select()
